# Refresh the cryptocurrency Price (D) and Volume(1h) (E) columns with the latest
# scraped snapshot (GitHub Actions run on 2024-09-29).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D ("Price") cells are stored as text (e.g. "65.747.10" uses dots as
# thousands separators, and values such as "1.00"/"8.00" rely on trailing zeros)
# so they must stay text after the update. Simply assigning a numeric-looking
# string to Range.Value lets Excel auto-convert it to a real number, which would
# drop the formatting/trailing zeros and change the cell type. Prefixing the new
# value with a literal leading apostrophe forces Excel to keep it as text; setting
# Style back to "Normal" afterwards clears the quote-prefix formatting Excel
# applies for that apostrophe so the cell keeps its original (unstyled) look.

$ws.Range("D2").Value = '''65.747.10'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.41%  '
$ws.Range("D3").Value = '''2.671.89'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.81%  '
$ws.Range("D5").Value = '''598.68'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -2.04%  '
$ws.Range("D6").Value = '''157.36'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.72%  '
$ws.Range("E7").Value = '  -0.01%  '
$ws.Range("E8").Value = '  +4.23%  '
$ws.Range("E9").Value = '  +2.22%  '
$ws.Range("E10").Value = '  -0.97%  '
$ws.Range("D11").Value = '''5.81'
$ws.Range("D11").Style = "Normal"
$ws.Range("E12").Value = '  -0.16%  '
$ws.Range("D13").Value = '''29.04'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -3.53%  '
$ws.Range("D14").Value = '''0.0000199'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -4.18%  '
$ws.Range("D15").Value = '''3.149.64'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.90%  '
$ws.Range("D16").Value = '''65.629.43'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.39%  '
$ws.Range("D17").Value = '''2.672.69'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.60%  '
$ws.Range("D18").Value = '''12.66'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.99%  '
$ws.Range("D19").Value = '''4.79'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -2.22%  '
$ws.Range("D20").Value = '''351.94'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.76%  '
$ws.Range("D21").Value = '''7.50'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -3.33%  '
$ws.Range("E22").Value = '  +0.06%  '
$ws.Range("D23").Value = '''69.18'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -3.05%  '
$ws.Range("E24").Value = '  -0.48%  '
$ws.Range("D25").Value = '''9.73'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -2.02%  '
$ws.Range("E26").Value = '  +2.67%  '
$ws.Range("D27").Value = '''1.60'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -3.82%  '
$ws.Range("E28").Value = '  -3.75%  '
$ws.Range("D29").Value = '''8.00'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -3.56%  '
$ws.Range("D31").Value = '''534.60'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.29%  '
$ws.Range("E32").Value = '  -3.67%  '
$ws.Range("E33").Value = '  -0.76%  '
$ws.Range("D34").Value = '''6.47'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -3.27%  '
$ws.Range("D35").Value = '''5.48'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.11%  '
$ws.Range("D36").Value = '''0.423'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -2.49%  '
$ws.Range("D37").Value = '''20.57'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.90%  '
$ws.Range("E38").Value = '  -0.03%  '
$ws.Range("D39").Value = '''156.98'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -3.51%  '
$ws.Range("E40").Value = '  -2.91%  '
$ws.Range("D41").Value = '''1.00'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.03%  '
$ws.Range("D42").Value = '''163.03'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -3.03%  '
$ws.Range("E43").Value = '  -1.51%  '
$ws.Range("D44").Value = '''2.34'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +1.32%  '
$ws.Range("D45").Value = '''0.0611'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -4.01%  '
$ws.Range("D46").Value = '''22.69'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -4.78%  '
$ws.Range("D47").Value = '''0.640'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -2.80%  '
$ws.Range("D48").Value = '''0.0258'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -3.55%  '
$ws.Range("D49").Value = '''0.0₆0253'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +7.02%  '
$ws.Range("E50").Value = '  -1.05%  '
$ws.Range("D51").Value = '''19.94'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -4.62%  '
